$d = $word.ActiveDocument

# Locate the unique paragraph: "PCount:S" and "PCount:M" instead of "M3PS" and "M3PM",
foreach ($para in $d.Paragraphs) {
    $r = $para.Range
    if ($r.Text -like '*PCount:S*' -and $r.Text -like '*PCount:M*' -and $r.Text -like '*M3PS*') {
        $r.Font.StrikeThrough = $true
    }
}
